# Apply the "dialogic" sheet correlation-matrix update:
#  - insert a new "sentence_complexity" row/column right after "word_length"
#  - append a new "discourse_markers_d" row/column at the end
# Resulting matrix grows from 6x6 (A1:F6) to 8x8 (A1:H8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "dialogic"

# ------------------------------------------------------------------
# 1. Make room: insert a blank column at C and a blank row at 3 so the
#    new "sentence_complexity" series slots in right after word_length.
# ------------------------------------------------------------------
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("A3").EntireRow.Insert()

# ------------------------------------------------------------------
# 2. Grow the grid so there is space for the appended
#    "discourse_markers_d" row/column (column H, row 8).
# ------------------------------------------------------------------
$ws.Range("H1").Value = " "
$ws.Range("A8").Value = " "

# ------------------------------------------------------------------
# 3. Re-apply the header style (bold, centered, thin border - same
#    as the original header cells) to every label cell, since insert
#    operations and fresh cells don't automatically inherit it.
# ------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4. Column / row headers (labels)
# ------------------------------------------------------------------
$labels = @("word_length", "sentence_complexity", "personal_pronoun_d", "passive_voice_d", "interjection_d", "modal_verb_d", "discourse_markers_d")

for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $labels[$i]
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}

# ------------------------------------------------------------------
# 5. Correlation values (symmetric 7x7 matrix, rows/cols 2..8 <-> B..H)
# ------------------------------------------------------------------
$matrix = @(
    @(1, 0.2939859033294212, -0.06672474766691336, 0.4035961329072622, -0.3153996724127018, 0.3019432252019841, 0.2403625567609906),
    @(0.2939859033294212, 1, 0.7373785647927182, 0.5657222839412497, -0.1532390230953878, 0.4999688335256128, 0.7786104552505485),
    @(-0.06672474766691336, 0.7373785647927182, 1, 0.3795779975692216, -0.01415982120200794, 0.4135628921474821, 0.7332728826971066),
    @(0.4035961329072622, 0.5657222839412497, 0.3795779975692216, 1, -0.1227941226378256, 0.42186840413018, 0.5765897079294763),
    @(-0.3153996724127018, -0.1532390230953878, -0.01415982120200794, -0.1227941226378256, 1, -0.2609126291512586, 0.002302043467693686),
    @(0.3019432252019841, 0.4999688335256128, 0.4135628921474821, 0.42186840413018, -0.2609126291512586, 1, 0.4580457105838001),
    @(0.2403625567609906, 0.7786104552505485, 0.7332728826971066, 0.5765897079294763, 0.002302043467693686, 0.4580457105838001, 1)
)

for ($r = 0; $r -lt 7; $r++) {
    for ($c = 0; $c -lt 7; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $matrix[$r][$c]
    }
}

$ws.Range("A1").Select()

# Restore the originally-active sheet so we don't leave stray
# workbook/sheet "active tab" side effects behind.
$wb.Worksheets.Item(1).Activate()
